# iris128b-bom.xlsx manufacturing update
# - Designator (column B) text updated for the 0.1uF caps (row3), 100R resistors
#   (row5), 0-ohm jumpers (row6) and 0.01uF caps (row7) rows, reflecting the
#   final as-built reference designators (ranges instead of enumerated lists,
#   plus R1/R6 reassigned between the 100R and jumper groups).
# - Designator column widened so the (now longer) lists/ranges are readable,
#   and the row heights that were driven by the old wrapped text are
#   refreshed: the two rows whose text now fits on one line drop back to the
#   sheet's default height, the others are shortened to match their new
#   (shorter) wrapped content.
# - View refreshed: zoomed to 87% with the last-edited cell (E8) selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Designator (column B) updates ---------------------------------------
$ws.Range("B3").Value = "C1-C24"
$ws.Range("B5").Value = "R10, R11, R2, R3, R4, R5, R6, R7, R8, R9"
$ws.Range("B6").Value = "LVDS1, LVDS3, LVDS2, LVDS4, LVDS5, LVDS8, LVDS7, LVDS6, R1"
$ws.Range("B7").Value = "C25-32"

# --- Column B widened to fit the designator lists -------------------------
$ws.Columns.Item(2).ColumnWidth = 98.26

# --- Row heights refreshed for the new (shorter) wrapped text -------------
$ws.Rows.Item(2).RowHeight = 25
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).RowHeight = 25
$ws.Rows.Item(7).RowHeight = 50
$ws.Rows.Item(8).RowHeight = 25

# --- View: zoom + selection -------------------------------------------
[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 87
[void]$ws.Range("E8").Select()
